$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "Recorded By" list. Some rows record it as
# "dnasr281@gmail.com, System" - swap the order to "System, dnasr281@gmail.com"
# everywhere that exact value occurs, leaving any other values untouched.
$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$recordedByColumn = $ws.Columns.Item(7)
$used = $excel.Intersect($recordedByColumn, $ws.UsedRange)

foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
